$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-27 Friday" "2024-12-28 Saturday"

Replace-Text "677×9=6093" "910×6=5460"
Replace-Text "893×6=5358" "221×3=663"
Replace-Text "168×5=840" "569×6=3414"
Replace-Text "407×3=1221" "277×2=554"
Replace-Text "107×9=963" "235×6=1410"

Replace-Text "712×6=4272" "289×8=2312"
Replace-Text "101×4=404" "883×5=4415"
Replace-Text "411×3=1233" "526×5=2630"
Replace-Text "587×8=4696" "122×2=244"
Replace-Text "722×6=4332" "648×8=5184"

Replace-Text "644×8=5152" "690×3=2070"
Replace-Text "430×6=2580" "698×4=2792"
Replace-Text "292×4=1168" "390×2=780"
Replace-Text "823×7=5761" "827×4=3308"
Replace-Text "660×6=3960" "203×2=406"

Replace-Text "367×3=1101" "898×8=7184"
Replace-Text "262×8=2096" "406×5=2030"
Replace-Text "735×3=2205" "965×7=6755"
Replace-Text "342×8=2736" "820×2=1640"
Replace-Text "244×2=488" "951×9=8559"

Replace-Text "854×9=7686" "469×5=2345"
Replace-Text "770×8=6160" "977×6=5862"
Replace-Text "401×7=2807" "891×9=8019"
Replace-Text "939×3=2817" "765×3=2295"
Replace-Text "664×9=5976" "790×8=6320"
